# Updates cryptos list values (Price / Volume(1h) columns, and two swapped rows)
# to match the refreshed data pulled on Tue Jan 23 15:32:33 UTC 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to be written as text (even when the new value looks
# numeric, e.g. "0.999" or "86.72") so the cell keeps the original text-cell
# semantics, then strip the temporary text NumberFormat back off again.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range('D2') '39.037.53'
Set-TextValue $ws.Range('E2') '  -3.44%  '
# Row 3
Set-TextValue $ws.Range('D3') '2.214.09'
Set-TextValue $ws.Range('E3') '  -6.19%  '
# Row 4
Set-TextValue $ws.Range('E4') '  +0.07%  '
# Row 5
Set-TextValue $ws.Range('D5') '296.03'
Set-TextValue $ws.Range('E5') '  -4.16%  '
# Row 6
Set-TextValue $ws.Range('D6') '82.14'
Set-TextValue $ws.Range('E6') '  -4.72%  '
# Row 7
Set-TextValue $ws.Range('D7') '0.509'
Set-TextValue $ws.Range('E7') '  -2.92%  '
# Row 8
Set-TextValue $ws.Range('E8') '  +0.05%  '
# Row 9
Set-TextValue $ws.Range('D9') '0.468'
Set-TextValue $ws.Range('E9') '  -4.41%  '
# Row 10
Set-TextValue $ws.Range('E10') '  -7.32%  '
# Row 11
Set-TextValue $ws.Range('D11') '28.67'
Set-TextValue $ws.Range('E11') '  -5.45%  '
# Row 12
Set-TextValue $ws.Range('D12') '46.65'
Set-TextValue $ws.Range('E12') '  -11.04%  '
# Row 13
Set-TextValue $ws.Range('E13') '  -1.99%  '
# Row 14
Set-TextValue $ws.Range('E14') '  -6.34%  '
# Row 15
Set-TextValue $ws.Range('D15') '6.17'
Set-TextValue $ws.Range('E15') '  -4.77%  '
# Row 16
Set-TextValue $ws.Range('D16') '14.04'
Set-TextValue $ws.Range('E16') '  -5.61%  '
# Row 17
Set-TextValue $ws.Range('D17') '2.216.11'
Set-TextValue $ws.Range('E17') '  -6.10%  '
# Row 18
Set-TextValue $ws.Range('E18') '  -5.40%  '
# Row 19
Set-TextValue $ws.Range('D19') '38.941.10'
Set-TextValue $ws.Range('E19') '  -3.53%  '
# Row 20
Set-TextValue $ws.Range('D20') '0.0₃0866'
Set-TextValue $ws.Range('E20') '  -4.03%  '
# Row 21
Set-TextValue $ws.Range('D21') '5.69'
Set-TextValue $ws.Range('E21') '  -6.35%  '
# Row 22
Set-TextValue $ws.Range('D22') '64.60'
Set-TextValue $ws.Range('E22') '  -5.21%  '
# Row 23
Set-TextValue $ws.Range('D23') '10.12'
Set-TextValue $ws.Range('E23') '  -4.78%  '
# Row 24
Set-TextValue $ws.Range('D24') '226.02'
Set-TextValue $ws.Range('E24') '  -2.87%  '
# Row 25
Set-TextValue $ws.Range('E25') '  +0.10%  '
# Row 26
Set-TextValue $ws.Range('D26') '2.38'
Set-TextValue $ws.Range('E26') '  -7.38%  '
# Row 27
Set-TextValue $ws.Range('D27') '1.76'
Set-TextValue $ws.Range('E27') '  -1.66%  '
# Row 28
Set-TextValue $ws.Range('D28') '22.41'
Set-TextValue $ws.Range('E28') '  -4.20%  '
# Row 29
Set-TextValue $ws.Range('E29') '  -1.43%  '
# Row 30
Set-TextValue $ws.Range('D30') '9.03'
Set-TextValue $ws.Range('E30') '  -2.31%  '
# Row 31
Set-TextValue $ws.Range('B31') 'InjectiveProtocol'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D31') '31.70'
Set-TextValue $ws.Range('E31') '  -5.14%  '
# Row 32
Set-TextValue $ws.Range('B32') 'Monero'
Set-TextValue $ws.Range('C32') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D32') '147.52'
Set-TextValue $ws.Range('E32') '  -2.63%  '
# Row 33
Set-TextValue $ws.Range('D33') '0.999'
Set-TextValue $ws.Range('E33') '  -0.23%  '
# Row 34
Set-TextValue $ws.Range('D34') '4.81'
Set-TextValue $ws.Range('E34') '  -6.56%  '
# Row 35
Set-TextValue $ws.Range('D35') '0.0695'
Set-TextValue $ws.Range('E35') '  -3.75%  '
# Row 36
Set-TextValue $ws.Range('D36') '2.32'
Set-TextValue $ws.Range('E36') '  -4.85%  '
# Row 37
Set-TextValue $ws.Range('E37') '  -3.40%  '
# Row 38
Set-TextValue $ws.Range('D38') '2.64'
Set-TextValue $ws.Range('E38') '  -2.99%  '
# Row 39
Set-TextValue $ws.Range('D39') '0.0948'
# Row 40
Set-TextValue $ws.Range('D40') '14.76'
Set-TextValue $ws.Range('E40') '  -5.96%  '
# Row 41
Set-TextValue $ws.Range('D41') '1.61'
Set-TextValue $ws.Range('E41') '  -4.36%  '
# Row 42
Set-TextValue $ws.Range('D42') '3.69'
Set-TextValue $ws.Range('E42') '  -2.71%  '
# Row 43
Set-TextValue $ws.Range('D43') '1.904.01'
Set-TextValue $ws.Range('E43') '  -2.16%  '
# Row 44
Set-TextValue $ws.Range('D44') '0.0257'
Set-TextValue $ws.Range('E44') '  -3.82%  '
# Row 45
Set-TextValue $ws.Range('E45') '  -16.88%  '
# Row 46
Set-TextValue $ws.Range('D46') '9.00'
Set-TextValue $ws.Range('E46') '  -4.06%  '
# Row 47
Set-TextValue $ws.Range('D47') '16.04'
Set-TextValue $ws.Range('E47') '  -8.00%  '
# Row 48
Set-TextValue $ws.Range('D48') '2.60'
Set-TextValue $ws.Range('E48') '  -2.82%  '
# Row 49
Set-TextValue $ws.Range('D49') '2.426.36'
Set-TextValue $ws.Range('E49') '  -6.34%  '
# Row 50
Set-TextValue $ws.Range('D50') '70.43'
Set-TextValue $ws.Range('E50') '  -2.18%  '
# Row 51
Set-TextValue $ws.Range('D51') '86.72'
Set-TextValue $ws.Range('E51') '  -6.06%  '
